$d = $word.ActiveDocument
$d.Content.Find.Execute("2025-12-16 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-12-17 Wednesday", 2) | Out-Null

$tbl = $d.Tables.Item(1)
$tbl.Cell(1,1).Range.Text = "23÷7="
$tbl.Cell(1,2).Range.Text = "47÷2="
$tbl.Cell(1,3).Range.Text = "56÷8="
$tbl.Cell(1,4).Range.Text = "66÷9="
$tbl.Cell(1,5).Range.Text = "80÷2="
$tbl.Cell(5,1).Range.Text = "21÷4="
$tbl.Cell(5,2).Range.Text = "44÷5="
$tbl.Cell(5,3).Range.Text = "28÷8="
$tbl.Cell(5,4).Range.Text = "71÷9="
$tbl.Cell(5,5).Range.Text = "65÷9="
$tbl.Cell(9,1).Range.Text = "27÷5="
$tbl.Cell(9,2).Range.Text = "59÷2="
$tbl.Cell(9,3).Range.Text = "71÷8="
$tbl.Cell(9,4).Range.Text = "88÷6="
$tbl.Cell(9,5).Range.Text = "17÷9="
$tbl.Cell(13,1).Range.Text = "97÷6="
$tbl.Cell(13,2).Range.Text = "28÷8="
$tbl.Cell(13,3).Range.Text = "71÷2="
$tbl.Cell(13,4).Range.Text = "78÷7="
$tbl.Cell(13,5).Range.Text = "67÷5="
$tbl.Cell(17,1).Range.Text = "32÷8="
$tbl.Cell(17,2).Range.Text = "67÷3="
$tbl.Cell(17,3).Range.Text = "50÷5="
$tbl.Cell(17,4).Range.Text = "71÷9="
$tbl.Cell(17,5).Range.Text = "76÷4="
